$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive leading text.
$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $pText = $paras.Item($i).Range.Text
    if ($pText -like "When the store needs to buy new DVDs from a supplier*") {
        $target = $paras.Item($i)
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$r = $target.Range
$pStart = $r.Start
$pText = $r.Text

# Find the "store manager" that precedes " will find the supplier record"
# (there is another "store manager" later in the same paragraph).
$needle = "store manager will find the supplier record"
$localIdx = $pText.IndexOf($needle)
if ($localIdx -lt 0) {
    throw "Could not locate 'store manager will find the supplier record' text"
}

$hlStart = $pStart + $localIdx
$hlEnd = $hlStart + "store manager".Length

# Narrow the resolved paragraph range down to just "store manager" and
# apply yellow highlighting to it, splitting the run.
$r.SetRange($hlStart, $hlEnd)
$r.Font.HighlightColorIndex = 7
